$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Last status check on" timestamp text in F1
$ws.Range("F1").Value = "Last status check on: 02.02.2022 05:45"

# 2. Convert D5 from inline text "+0.6" to numeric 0.6
$ws.Range("D5").Value = 0.6

# 3. Convert E5 from inline text date to a real numeric date/time value,
#    matching the style used by the other rows (E6:E10) in column E.
$ws.Range("E5").Value = 44594.2311574074
$ws.Range("E5").NumberFormat = $ws.Range("E6").NumberFormat
